$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# LINO TUMBACO VICENTE JAVIER / PIEDRA SINTERIZADA
$ws1.Range("L14").Value = 565.25
# PAREDES ORTIZ MARIA INES / 240X80 PORCELANATO
$ws1.Range("D18").Value = 1900.8

# Row 29 summary counts ("N de 27") recalculated after the above changes
$ws1.Range("D29").Value = "1 de 27"
$ws1.Range("L29").Value = "2 de 27"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# LINO TUMBACO VICENTE JAVIER / julio
$ws2.Range("F14").Value = 6.93
# PAREDES ORTIZ MARIA INES / julio
$ws2.Range("F18").Value = 4078.08
# TOTAL row / julio
$ws2.Range("F29").Value = 6130.32

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 240X80 PORCELANATO
$ws3.Range("D3").Value = 1900.8
$ws3.Range("E3").Value = 1219.3145
$ws3.Range("F3").Value = 0.6092084120630828

# PIEDRA SINTERIZADA
$ws3.Range("D15").Value = 874.72
$ws3.Range("E15").Value = 763.28
$ws3.Range("F15").Value = 0.5340170940170941

# TOTAL
$ws3.Range("D19").Value = 6364.79
$ws3.Range("E19").Value = 20816.52093005039
$ws3.Range("F19").Value = 0.6092084120630828
